# "Added support for longer quotes, fixed surplus numnber"
#
# - The "surplus number" fudge-factor (1.0565) in column K is reset back to a
#   plain 1 for the first row of each size group (rows 16, 17, 20, 23, 26,
#   29, 32, 33, 37 on the "Customer Quote" sheet).
# - The active selection is moved from A37 to D8 (to support longer quotes /
#   more rows being visible near the top of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

$rowsToFix = @(16, 17, 20, 23, 26, 29, 32, 33, 37)
foreach ($r in $rowsToFix) {
    $ws.Cells.Item($r, 11).Value = 1   # column K = 11
}

$ws.Activate()
$ws.Range("D8").Select()
